# ADD results from server
$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 973.9537847600009
$ws.Range("E2").Value = 28982.37596598056
$ws.Range("I2").Value = 16175.28135478
$ws.Range("L2").Value = 48524.529503538
$ws.Range("M2").Value = 10590.587968015
$ws.Range("N2").Value = 7152.019986098924
$ws.Range("O2").Value = 6979.505869462286

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 5712.560177842886
$ws.Range("E2").Value = 56106.05588781912
$ws.Range("I2").Value = 44217.8984721661
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 21984.28023276101
$ws.Range("N2").Value = 10587.30466466927
$ws.Range("O2").Value = 12059.56090694894

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 2861.961401238371
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15122.94790489287
$ws.Range("O2").Value = 14757.33642912432

# Sheet "2040"
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 2861.961401238371
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15227.20816887006
$ws.Range("O2").Value = 14757.33642912432

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 6302.873118834019
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15761.83005432624
$ws.Range("O2").Value = 17090.83148627414

# Sheet "2050"
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 6302.873118834019
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15761.83005432624
$ws.Range("O2").Value = 17090.83148627414
